# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Maldivas" to sit alphabetically before "Islas Feroe" ---
# Row 131 used to be "Islas Feroe", row 132 "Ruanda", row 133 "Maldivas".
# New order: row 131 "Maldivas" (with refreshed stats), row 132 "Islas Feroe"
# (stats carried over unchanged), row 133 "Ruanda" (stats carried over unchanged).
$ws.Range("A131").Value = "Maldivas"
$ws.Range("B131").Value = 191
$ws.Range("C131").Value = 14
$ws.Range("D131").Value = 17
$ws.Range("E131").Value = 174
$ws.Range("F131").Value = 2

$ws.Range("A132").Value = "Islas Feroe"
$ws.Range("B132").Value = 187
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 178
$ws.Range("E132").Value = 9

$ws.Range("A133").Value = "Ruanda"
$ws.Range("B133").Value = 183
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 88
$ws.Range("E133").Value = 95
$ws.Range("F133").Value = 0

# --- Refresh case numbers for several countries ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 961500
$ws.Range("C4").Value = 849
$ws.Range("E4").Value = 789061
$ws.Range("G4").Value = 21
$ws.Range("H4").Value = 54277

# Alemania (row 8)
$ws.Range("B8").Value = 156727
$ws.Range("C8").Value = 214
$ws.Range("E8").Value = 41047
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 5880

# Argentina (row 56)
$ws.Range("D56").Value = 1107
$ws.Range("E56").Value = 2486
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 187

# Bulgaria (row 84)
$ws.Range("B84").Value = 1300
$ws.Range("C84").Value = 53
$ws.Range("E84").Value = 1039
$ws.Range("F84").Value = 41
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 56

# Estado de Palestina (row 118)
$ws.Range("D118").Value = 83
$ws.Range("E118").Value = 257

# Barbados (row 154)
$ws.Range("D154").Value = 32
$ws.Range("E154").Value = 41
